$d = $word.ActiveDocument

function Add-EmptyParagraph($afterParagraph) {
    $afterParagraph.Range.InsertParagraphAfter() | Out-Null
    return $d.Paragraphs.Last
}

function Add-TextParagraph($afterParagraph, [string]$text) {
    $p = Add-EmptyParagraph $afterParagraph
    $p.Range.InsertAfter($text) | Out-Null
    return $p
}

# Start from the last paragraph in the document (the existing final entry).
$cur = $d.Paragraphs.Last

# Two blank paragraphs.
$cur = Add-EmptyParagraph $cur
$cur = Add-EmptyParagraph $cur

# A paragraph that is just a manual page break.
# Building it by splitting a two-character placeholder run around the break
# (mirrors how Word itself splits a run when a break lands inside it) keeps
# the paragraph/run XML shape consistent with a normal in-place page break.
$cur = Add-EmptyParagraph $cur
$cur.Range.InsertAfter("AB") | Out-Null
$breakPos = $cur.Range.Start + 1
$d.Range($breakPos, $breakPos).InsertBreak(7) | Out-Null
$bRange = $d.Range($cur.Range.Start + 2, $cur.Range.Start + 3)
$bRange.Delete() | Out-Null
$aRange = $d.Range($cur.Range.Start, $cur.Range.Start + 1)
$aRange.Delete() | Out-Null

# 03/08/22 stand up entries.
$cur = Add-TextParagraph $cur "03/08/22"

$cur = Add-TextParagraph $cur "Devin"
$cur = Add-TextParagraph $cur "Yesterday worked on risk assessment and starting our front end by creating a react app. Today plan on finishing at least wireframe and the home page so that we can send Morgan our wire frame and have a good start for the website.Blockers - lack of experience with jsx."
$cur = Add-EmptyParagraph $cur

$cur = Add-TextParagraph $cur "Hodan"
$cur = Add-TextParagraph $cur "Yesterday worked on risk assessment and starting frontend and backend and planning where to go from there. Today working on home page, footer, and wireframe. Blockers, trying to understand coding with jsx."
$cur = Add-EmptyParagraph $cur

$cur = Add-TextParagraph $cur "Toseef"
$cur = Add-TextParagraph $cur "Yesterday worked on risk assessment, planning structure for frontend. Today working on home page, header/navbar, and wireframe. Blockers, lack of experience with react/jsx."
$cur = Add-EmptyParagraph $cur

$cur = Add-TextParagraph $cur "Waseem"
$cur = Add-TextParagraph $cur "Unable to attend (has covid which may be the cause?)"

Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
